$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 59.8125
$ws.Range("I11").Value = 59.8125
$ws.Range("K11").Value = 59.8125
$ws.Range("M11").Value = 80.1875

$ws.Range("H17").Value = 1899.0286
$ws.Range("J17").Value = 1899.0286
$ws.Range("L17").Value = 5697.085800000001
$ws.Range("N17").Value = -6033.085800000001

$ws.Range("H69").Value = 8456.286
$ws.Range("I69").Value = 5053
$ws.Range("K69").Value = 15159
$ws.Range("M69").Value = -14285

$ws.Range("H72").Value = 8456.286
$ws.Range("I72").Value = 5053
$ws.Range("K72").Value = 45477
$ws.Range("M72").Value = -41109

$ws.Range("H74").Value = 7449.75
$ws.Range("I74").Value = 4299.875
$ws.Range("K74").Value = 4299.875
$ws.Range("M74").Value = -3363.875

$ws.Range("H77").Value = 7449.75
$ws.Range("I77").Value = 4299.875
$ws.Range("K77").Value = 21499.375
$ws.Range("M77").Value = -16819.375

$ws.Range("H94").Value = 2961.4443
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H116").Value = 40382.57
$ws.Range("J116").Value = 9987.299999999999
$ws.Range("L116").Value = 9987.299999999999
$ws.Range("N116").Value = -16871.3

$ws.Range("H125").Value = 2307.75
$ws.Range("I125").Value = 1964.875
$ws.Range("K125").Value = 17683.875
$ws.Range("M125").Value = -15223.875

$ws.Range("H132").Value = 2412.5862
$ws.Range("I132").Value = 2355.9285
$ws.Range("K132").Value = 7067.7855
$ws.Range("M132").Value = -4537.7855

$ws.Range("H137").Value = 2209.3438
$ws.Range("I137").Value = 1264.2858
$ws.Range("K137").Value = 3792.8574
$ws.Range("M137").Value = -1242.8574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4380.625
$ws.Range("I32").Value = 4657.925
$ws.Range("K32").Value = 4657.925
$ws.Range("M32").Value = -4370.925

$ws.Range("H74").Value = 2275.3489
$ws.Range("I74").Value = 1508.1471
$ws.Range("K74").Value = 1508.1471
$ws.Range("M74").Value = -634.1470999999999

$ws.Range("H77").Value = 2275.3489
$ws.Range("I77").Value = 1508.1471
$ws.Range("K77").Value = 7540.7355
$ws.Range("M77").Value = -3172.7355

$ws.Range("H94").Value = 42500
$ws.Range("J94").Value = 42500
$ws.Range("L94").Value = 42500
$ws.Range("N94").Value = -44302

$ws.Range("H122").Value = 4224.6665
$ws.Range("I122").Value = 3945.6
$ws.Range("K122").Value = 11836.8
$ws.Range("M122").Value = -9386.799999999999

$ws.Range("H132").Value = 3207.0886
$ws.Range("I132").Value = 1860.8966
$ws.Range("K132").Value = 5582.6898
$ws.Range("M132").Value = -3052.6898

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1585.1852
$ws.Range("I105").Value = 1445.8636
$ws.Range("K105").Value = 1445.8636
$ws.Range("M105").Value = 301.1364000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1528.0555
$ws.Range("I22").Value = 763
$ws.Range("J22").Value = 3517.2
$ws.Range("K22").Value = 763
$ws.Range("L22").Value = 3517.2
$ws.Range("M22").Value = -413
$ws.Range("N22").Value = -4217.2

$ws.Range("H31").Value = 4656
$ws.Range("I31").Value = 3186.2632
$ws.Range("K31").Value = 3186.2632
$ws.Range("M31").Value = -2891.2632

$ws.Range("H34").Value = 4656
$ws.Range("I34").Value = 3186.2632
$ws.Range("K34").Value = 3186.2632
$ws.Range("M34").Value = -2984.2632

$ws.Range("H132").Value = 1722.1428
$ws.Range("I132").Value = 1741.1538
$ws.Range("K132").Value = 5223.4614
$ws.Range("M132").Value = -2693.4614

$ws.Range("H134").Value = 1962.7037
$ws.Range("I134").Value = 1949.0454
$ws.Range("K134").Value = 5847.1362
$ws.Range("M134").Value = -3312.1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 510.4
$ws.Range("J98").Value = 602.8
$ws.Range("L98").Value = 1808.4
$ws.Range("N98").Value = -4804.4

$ws.Range("H114").Value = 4260.8184
$ws.Range("I114").Value = 311.25
$ws.Range("J114").Value = 6517.7144
$ws.Range("K114").Value = 933.75
$ws.Range("L114").Value = 19553.1432
$ws.Range("M114").Value = 2320.25
$ws.Range("N114").Value = -26061.1432

$ws.Range("H128").Value = 260000
$ws.Range("I128").Value = 260000
$ws.Range("K128").Value = 780000
$ws.Range("M128").Value = -775020

$ws.Range("H137").Value = 4516.7144
$ws.Range("J137").Value = 6590.5
$ws.Range("L137").Value = 19771.5
$ws.Range("N137").Value = -29971.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 62587.5
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

$ws.Range("H97").Value = 1516.1842
$ws.Range("I97").Value = 596.7826
$ws.Range("J97").Value = 2925.9333
$ws.Range("K97").Value = 596.7826
$ws.Range("L97").Value = 2925.9333
$ws.Range("M97").Value = -100.7826
$ws.Range("N97").Value = -3917.9333

$ws.Range("H102").Value = 2014.1
$ws.Range("I102").Value = 1267.6666
$ws.Range("J102").Value = 4999.8335
$ws.Range("K102").Value = 1267.6666
$ws.Range("L102").Value = 4999.8335
$ws.Range("M102").Value = 354.3334
$ws.Range("N102").Value = -8243.833500000001

$ws.Range("H113").Value = 2292.0715
$ws.Range("I113").Value = 1207.7
$ws.Range("K113").Value = 1207.7
$ws.Range("M113").Value = 962.3

$ws.Range("H122").Value = 4363.05
$ws.Range("I122").Value = 3488.4666
$ws.Range("K122").Value = 10465.3998
$ws.Range("M122").Value = -8015.399800000001

$ws.Range("H126").Value = 4164.8
$ws.Range("I126").Value = 3099
$ws.Range("J126").Value = 4431.25
$ws.Range("K126").Value = 9297
$ws.Range("L126").Value = 13293.75
$ws.Range("M126").Value = -6827
$ws.Range("N126").Value = -18233.75

$ws.Range("H132").Value = 1351.5625
$ws.Range("I132").Value = 1175.0834
$ws.Range("K132").Value = 3525.2502
$ws.Range("M132").Value = -995.2501999999999

$ws.Range("H135").Value = 95000
$ws.Range("J135").Value = 95000
$ws.Range("L135").Value = 95000
$ws.Range("N135").Value = -105140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10597.2
$ws.Range("I16").Value = 13079.75
$ws.Range("J16").Value = 667
$ws.Range("K16").Value = 13079.75
$ws.Range("L16").Value = 667
$ws.Range("M16").Value = -12909.75
$ws.Range("N16").Value = -1007

$ws.Range("H22").Value = 144941.72
$ws.Range("I22").Value = 250650.5
$ws.Range("K22").Value = 250650.5
$ws.Range("M22").Value = -250355.5

$ws.Range("H27").Value = 144941.72
$ws.Range("I27").Value = 250650.5
$ws.Range("K27").Value = 250650.5
$ws.Range("M27").Value = -250543.5

$ws.Range("H55").Value = 980.96295
$ws.Range("I55").Value = 240.9375
$ws.Range("K55").Value = 240.9375
$ws.Range("M55").Value = -67.9375

$ws.Range("H122").Value = 3370.311
$ws.Range("I122").Value = 3145.3242
$ws.Range("K122").Value = 9435.972600000001
$ws.Range("M122").Value = -6985.972600000001

$ws.Range("H136").Value = 4608.6665
$ws.Range("I136").Value = 4212.1313
$ws.Range("J136").Value = 8375.75
$ws.Range("K136").Value = 12636.3939
$ws.Range("L136").Value = 25127.25
$ws.Range("M136").Value = -10086.3939
$ws.Range("N136").Value = -30227.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4147.857
$ws.Range("I126").Value = 4037.9678
$ws.Range("K126").Value = 12113.9034
$ws.Range("M126").Value = -9643.903399999999

$ws.Range("H132").Value = 5735.1714
$ws.Range("I132").Value = 4861.207
$ws.Range("K132").Value = 14583.621
$ws.Range("M132").Value = -12053.621
